$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Split the single "User Name" column into "First Name" and "Last Name".
# Insert a new column at C (pushing the old C..T columns one to the right -
# "Password" ends up moving from T to U automatically), then fill in the
# header + data for the two new name columns.
# ---------------------------------------------------------------------------
$ws.Range("C1").EntireColumn.Insert()

$ws.Range("B1").Value = "First Name"
$ws.Range("C1").Value = "Last Name"

$ws.Range("B2").Value = "Test"
$ws.Range("C2").Value = "User"
